$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$textCells = @("D4", "D5", "D7", "D8", "D10", "D11", "D13", "D14", "D16", "D18", "D19", "D22", "D23", "D25", "D30", "D31", "D32", "D38", "D39", "D40", "D42", "D44", "D46", "D47", "D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "30.768.54"
$ws.Range("E2").Value = "  +2.54%  "
$ws.Range("D3").Value = "1.694.58"
$ws.Range("E3").Value = "  +3.49%  "
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.22%  "
$ws.Range("D5").Value = "221.85"
$ws.Range("E5").Value = "  +2.88%  "
$ws.Range("E6").Value = "  +0.26%  "
$ws.Range("D7").Value = "0.997"
$ws.Range("E7").Value = "  -0.21%  "
$ws.Range("D8").Value = "30.96"
$ws.Range("E8").Value = "  +4.06%  "
$ws.Range("E9").Value = "  +2.25%  "
$ws.Range("D10").Value = "0.0627"
$ws.Range("E10").Value = "  +2.10%  "
$ws.Range("D11").Value = "0.0905"
$ws.Range("E11").Value = "  -1.29%  "
$ws.Range("D12").Value = "1.940.74"
$ws.Range("E12").Value = "  +3.69%  "
$ws.Range("D13").Value = "10.71"
$ws.Range("E13").Value = "  +11.26%  "
$ws.Range("D14").Value = "0.623"
$ws.Range("E14").Value = "  +7.63%  "
$ws.Range("D15").Value = "1.697.50"
$ws.Range("E15").Value = "  +3.66%  "
$ws.Range("D16").Value = "4.02"
$ws.Range("E16").Value = "  +2.71%  "
$ws.Range("D17").Value = "30.798.59"
$ws.Range("D18").Value = "66.48"
$ws.Range("E18").Value = "  +2.38%  "
$ws.Range("D19").Value = "250.07"
$ws.Range("E19").Value = "  +0.33%  "
$ws.Range("E20").Value = "  +1.51%  "
$ws.Range("E21").Value = "  -0.27%  "
$ws.Range("D22").Value = "10.25"
$ws.Range("E22").Value = "  +4.88%  "
$ws.Range("D23").Value = "4.30"
$ws.Range("E23").Value = "  +2.39%  "
$ws.Range("E24").Value = "  +2.03%  "
$ws.Range("D25").Value = "157.51"
$ws.Range("E25").Value = "  -1.17%  "
$ws.Range("E26").Value = "  +1.47%  "
$ws.Range("E27").Value = "  +0.26%  "
$ws.Range("E28").Value = "  +1.33%  "
$ws.Range("E29").Value = "  -0.29%  "
$ws.Range("D30").Value = "0.0501"
$ws.Range("E30").Value = "  +1.94%  "
$ws.Range("D31").Value = "1.15"
$ws.Range("E31").Value = "  +1.34%  "
$ws.Range("D32").Value = "3.49"
$ws.Range("E32").Value = "  +2.76%  "
$ws.Range("D33").Value = "1.520.16"
$ws.Range("E33").Value = "  +5.87%  "
$ws.Range("E34").Value = "  +2.95%  "
$ws.Range("E35").Value = "  +4.85%  "
$ws.Range("E36").Value = "  -0.41%  "
$ws.Range("E37").Value = "  +4.75%  "
$ws.Range("D38").Value = "80.71"
$ws.Range("E38").Value = "  +4.40%  "
$ws.Range("B39").Value = "ImmutableX"
$ws.Range("C39").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D39").Value = "0.587"
$ws.Range("E39").Value = "  +4.95%  "
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").Value = "2.71"
$ws.Range("E40").Value = "  -5.32%  "
$ws.Range("E41").Value = "  +1.16%  "
$ws.Range("D42").Value = "0.857"
$ws.Range("E42").Value = "  +2.08%  "
$ws.Range("E43").Value = "  +1.61%  "
$ws.Range("D44").Value = "0.0504"
$ws.Range("E44").Value = "  +0.86%  "
$ws.Range("E45").Value = "  -1.52%  "
$ws.Range("B46").Value = "BitcoinSV"
$ws.Range("C46").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D46").Value = "53.02"
$ws.Range("E46").Value = "  -3.67%  "
$ws.Range("B47").Value = "PaxDollar"
$ws.Range("C47").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D47").Value = "0.996"
$ws.Range("E47").Value = "  -0.27%  "
$ws.Range("D48").Value = "1.828.90"
$ws.Range("E48").Value = "  +2.80%  "
$ws.Range("E49").Value = "  +0.22%  "
$ws.Range("D50").Value = "95.92"
$ws.Range("E50").Value = "  +6.02%  "
$ws.Range("D51").Value = "0.0₆0114"
$ws.Range("E51").Value = "  +3.21%  "
